$wb = $excel.ActiveWorkbook

$sheetNames = @("Preços Amazon", "Preços Shopee")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Copy the formatting of the existing "Lucro" header cell (F1) into the
    # three new header cells, then set their text.
    $ws.Range("F1").Copy($ws.Range("G1"))
    $ws.Range("F1").Copy($ws.Range("H1"))
    $ws.Range("F1").Copy($ws.Range("I1"))

    $ws.Range("G1").Value = "Lucro 3un."
    $ws.Range("H1").Value = "Lucro 5un."
    $ws.Range("I1").Value = "Lucro 10un."

    for ($row = 2; $row -le 5; $row++) {
        $lucro = $ws.Cells.Item($row, 6).Value2
        $ws.Cells.Item($row, 7).Value = $lucro * 3
        $ws.Cells.Item($row, 8).Value = $lucro * 5
        $ws.Cells.Item($row, 9).Value = $lucro * 10
    }
}
